$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1926.4286
$ws.Range("J17").Value = 2097.2727
$ws.Range("L17").Value = 6291.8181
$ws.Range("N17").Value = -6627.8181
$ws.Range("H62").Value = 10423992
$ws.Range("I62").Value = 12828376
$ws.Range("J62").Value = 4990.6665
$ws.Range("K62").Value = 12828376
$ws.Range("L62").Value = 4990.6665
$ws.Range("M62").Value = -12827752
$ws.Range("N62").Value = -6238.6665
$ws.Range("H65").Value = 10423992
$ws.Range("I65").Value = 12828376
$ws.Range("J65").Value = 4990.6665
$ws.Range("K65").Value = 64141880
$ws.Range("L65").Value = 24953.3325
$ws.Range("M65").Value = -64138760
$ws.Range("N65").Value = -31193.3325
$ws.Range("H70").Value = 1933.75
$ws.Range("I70").Value = 1814
$ws.Range("J70").Value = 2293
$ws.Range("K70").Value = 5442
$ws.Range("L70").Value = 6879
$ws.Range("M70").Value = -5172
$ws.Range("N70").Value = -7419
$ws.Range("H73").Value = 1933.75
$ws.Range("I73").Value = 1814
$ws.Range("J73").Value = 2293
$ws.Range("K73").Value = 5442
$ws.Range("L73").Value = 6879
$ws.Range("M73").Value = -4506
$ws.Range("N73").Value = -8751
$ws.Range("H93").Value = 42172
$ws.Range("J93").Value = 42172
$ws.Range("L93").Value = 42172
$ws.Range("N93").Value = -47164
$ws.Range("H115").Value = 549.8333
$ws.Range("J115").Value = 2000
$ws.Range("L115").Value = 6000
$ws.Range("N115").Value = -9134
$ws.Range("H116").Value = 18522168
$ws.Range("I116").Value = 22224600
$ws.Range("J116").Value = 9999
$ws.Range("K116").Value = 22224600
$ws.Range("L116").Value = 9999
$ws.Range("M116").Value = -22221158
$ws.Range("N116").Value = -16883
$ws.Range("H118").Value = 480.83334
$ws.Range("I118").Value = 190.11111
$ws.Range("K118").Value = 570.3333299999999
$ws.Range("M118").Value = 1086.66667
$ws.Range("H137").Value = 3926.1936
$ws.Range("I137").Value = 4736
$ws.Range("J137").Value = 3259.2942
$ws.Range("K137").Value = 14208
$ws.Range("L137").Value = 9777.882599999999
$ws.Range("M137").Value = -11658
$ws.Range("N137").Value = -14877.8826
$ws.Range("H138").Value = 3545.3137
$ws.Range("I138").Value = 2097.4666
$ws.Range("J138").Value = 4148.5835
$ws.Range("K138").Value = 6292.399800000001
$ws.Range("L138").Value = 12445.7505
$ws.Range("M138").Value = -1152.399800000001
$ws.Range("N138").Value = -22725.7505
$ws.Range("H141").Value = 2602.5
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").ClearContents()
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3009039.8
$ws.Range("I32").Value = 5299.25
$ws.Range("K32").Value = 5299.25
$ws.Range("M32").Value = -5012.25
$ws.Range("H45").Value = 1173.8182
$ws.Range("J45").Value = 500
$ws.Range("L45").Value = 500
$ws.Range("N45").Value = -1254
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9811334
$ws.Range("J31").Value = 4665.8647
$ws.Range("L31").Value = 4665.8647
$ws.Range("N31").Value = -5255.8647
$ws.Range("H34").Value = 9811334
$ws.Range("J34").Value = 4665.8647
$ws.Range("L34").Value = 4665.8647
$ws.Range("N34").Value = -5069.8647
$ws.Range("H141").Value = 610195.2
$ws.Range("J141").Value = 751494
$ws.Range("L141").Value = 751494
$ws.Range("N141").Value = -761854
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 31798.625
$ws.Range("I22").Value = 495
$ws.Range("J22").Value = 42233.168
$ws.Range("K22").Value = 1485
$ws.Range("L22").Value = 126699.504
$ws.Range("M22").Value = -1316
$ws.Range("N22").Value = -127037.504
$ws.Range("H26").Value = 497.4
$ws.Range("I26").Value = 497.4
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 1492.2
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -1204.2
$ws.Range("N26").ClearContents()
$ws.Range("H27").Value = 31798.625
$ws.Range("I27").Value = 495
$ws.Range("J27").Value = 42233.168
$ws.Range("K27").Value = 1485
$ws.Range("L27").Value = 126699.504
$ws.Range("M27").Value = -1383
$ws.Range("N27").Value = -126903.504
$ws.Range("H46").Value = 3449.5
$ws.Range("J46").Value = 3449.5
$ws.Range("L46").Value = 10348.5
$ws.Range("N46").Value = -10530.5
$ws.Range("H110").Value = 22782
$ws.Range("I110").Value = 17007.6
$ws.Range("K110").Value = 51022.8
$ws.Range("M110").Value = -46932.8
$ws.Range("H132").Value = 3217.4
$ws.Range("J132").Value = 4297.143
$ws.Range("L132").Value = 38674.287
$ws.Range("N132").Value = -43734.287
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6758.6665
$ws.Range("I80").Value = 3051.5
$ws.Range("J80").Value = 8612.25
$ws.Range("K80").Value = 3051.5
$ws.Range("L80").Value = 8612.25
$ws.Range("M80").Value = -2053.5
$ws.Range("N80").Value = -10608.25
$ws.Range("H83").Value = 6758.6665
$ws.Range("I83").Value = 3051.5
$ws.Range("J83").Value = 8612.25
$ws.Range("K83").Value = 15257.5
$ws.Range("L83").Value = 43061.25
$ws.Range("M83").Value = -10265.5
$ws.Range("N83").Value = -53045.25
$ws.Range("H126").Value = 15633952
$ws.Range("I126").Value = 21744158
$ws.Range("K126").Value = 65232474
$ws.Range("M126").Value = -65230004
$ws.Range("H132").Value = 3630.1968
$ws.Range("I132").Value = 3735.06
$ws.Range("J132").Value = 3153.5454
$ws.Range("K132").Value = 11205.18
$ws.Range("L132").Value = 9460.636200000001
$ws.Range("M132").Value = -8675.18
$ws.Range("N132").Value = -14520.6362
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1970.875
$ws.Range("I16").Value = 785.6
$ws.Range("K16").Value = 785.6
$ws.Range("M16").Value = -615.6
$ws.Range("H40").Value = 7499.5
$ws.Range("I40").Value = 5000
$ws.Range("K40").Value = 5000
$ws.Range("M40").Value = -4864
$ws.Range("H122").Value = 1819926.5
$ws.Range("I122").Value = 3329865.8
$ws.Range("J122").Value = 7999.6
$ws.Range("K122").Value = 9989597.399999999
$ws.Range("L122").Value = 23998.8
$ws.Range("M122").Value = -9987147.399999999
$ws.Range("N122").Value = -28898.8
$ws.Range("H132").Value = 4879
$ws.Range("I132").Value = 3352.5
$ws.Range("J132").Value = 7550.375
$ws.Range("K132").Value = 10057.5
$ws.Range("L132").Value = 22651.125
$ws.Range("M132").Value = -7527.5
$ws.Range("N132").Value = -27711.125
$ws.Range("H136").Value = 7156.6665
$ws.Range("I136").Value = 6409.2104
$ws.Range("K136").Value = 19227.6312
$ws.Range("M136").Value = -16677.6312
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 19080.438
$ws.Range("I51").Value = 22715.857
$ws.Range("J51").Value = 16252.889
$ws.Range("K51").Value = 22715.857
$ws.Range("L51").Value = 16252.889
$ws.Range("M51").Value = -22205.857
$ws.Range("N51").Value = -17272.889
